$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the (previously empty) row 3 with the new "$startTime" countdown entry ---
$ws.Range("A3").Value = '$startTime'
$ws.Range("B3").Value = 30

# Row 3 needs the same cell style as the other header/meta rows (row 1/2) -
# copy formats only so the values just set above are preserved.
$ws.Range("A2:B2").Copy()
$ws.Range("A3:B3").PasteSpecial(-4122)

# --- Update the first wave entry (row 5) in place ---
$ws.Range("A5").Value = 20
$ws.Range("B5").Value = 2

# --- Append the rest of the wave countdown / saucer-count table (rows 6-16) ---
$ws.Range("A6").Value = 15
$ws.Range("B6").Value = 5

$ws.Range("A7").Value = 10
$ws.Range("B7").Value = 6

$ws.Range("A8").Value = 10
$ws.Range("B8").Value = 7

$ws.Range("A9").Value = 10
$ws.Range("B9").Value = 8

$ws.Range("A10").Value = 10
$ws.Range("B10").Value = 9

$ws.Range("A11").Value = 8
$ws.Range("B11").Value = 10

$ws.Range("A12").Value = 8
$ws.Range("B12").Value = 12

$ws.Range("A13").Value = 8
$ws.Range("B13").Value = 14

$ws.Range("A14").Value = 8
$ws.Range("B14").Value = 16

$ws.Range("A15").Value = 6
$ws.Range("B15").Value = 18

$ws.Range("A16").Value = 6
$ws.Range("B16").Value = 20

# Column B on rows 7-16 carries the same style as B5/B4 above it (row 6 is left default).
$ws.Range("B5").Copy()
$ws.Range("B7:B16").PasteSpecial(-4122)

# New rows use the same explicit row height as the rest of the sheet.
$ws.Range($ws.Rows.Item(3), $ws.Rows.Item(16)).RowHeight = 15.75

# --- Restore the selection to match the authored workbook ---
$ws.Range("B5").Select()
